# Add a "Save" column (H) to the s_vals sheet, matching the style of the
# existing header row (e.g. column G's "sum" header).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style (bold font, border, centered alignment) from G1 to H1.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# Set the new header label and the data value.
$ws.Range("H1").Value = "Save"
$ws.Range("H2").Value = 1
